$d = $word.ActiveDocument
$count = 0

$rng = $d.Content
$ok = $rng.Find.Execute("英语", $false, $false, $false, $false, $false, $true, 1, $false, "英文", 1)
if (-not $ok) { Write-Output "MISS 0: 英语" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("葡萄牙语", $false, $false, $false, $false, $false, $true, 1, $false, "葡萄牙文", 1)
if (-not $ok) { Write-Output "MISS 1: 葡萄牙语" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("法语", $false, $false, $false, $false, $false, $true, 1, $false, "法文", 1)
if (-not $ok) { Write-Output "MISS 2: 法语" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("泰语", $false, $false, $false, $false, $false, $true, 1, $false, "泰文", 1)
if (-not $ok) { Write-Output "MISS 3: 泰语" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("越南语", $false, $false, $false, $false, $false, $true, 1, $false, "越南文", 1)
if (-not $ok) { Write-Output "MISS 4: 越南语" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("西班牙语", $false, $false, $false, $false, $false, $true, 1, $false, "西班牙文", 1)
if (-not $ok) { Write-Output "MISS 5: 西班牙语" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("英语", $false, $false, $false, $false, $false, $true, 1, $false, "英文", 1)
if (-not $ok) { Write-Output "MISS 6: 英语" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("简要", $false, $false, $false, $false, $false, $true, 1, $false, "簡介", 1)
if (-not $ok) { Write-Output "MISS 7: 简要" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("向目标国家的合作伙伴发送邀请他们参加为期一天的研讨会的电子邮件。 它将通过 customer.io 发送", $false, $false, $false, $false, $false, $true, 1, $false, "一封發送給目標國家夥伴的電子郵件，邀請他們參加為期一天的研討會。 將通過 customer.io 發送", 1)
if (-not $ok) { Write-Output "MISS 8: 向目标国家的合作伙伴发送邀请他们参加为期一天的研讨会的电子邮件。 它将通过 customer.io 发送" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("目标受众", $false, $false, $false, $false, $false, $true, 1, $false, "目標受眾", 1)
if (-not $ok) { Write-Output "MISS 9: 目标受众" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("目标国家的合作伙伴", $false, $false, $false, $false, $false, $true, 1, $false, "目標國家的夥伴", 1)
if (-not $ok) { Write-Output "MISS 10: 目标国家的合作伙伴" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute(": 在 [CITY] 与我们的团队会面 | [DATE] ", $false, $false, $false, $false, $false, $true, 1, $false, ": 在 [CITY] 與我們的團隊會面 | [DATE] ", 1)
if (-not $ok) { Write-Output "MISS 11: : 在 [CITY] 与我们的团队会面 | [DATE] " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("邀请参加 Deriv 合作伙伴研讨会", $false, $false, $false, $false, $false, $true, 1, $false, "邀請參加 Deriv 合作夥伴研討會", 1)
if (-not $ok) { Write-Output "MISS 12: 邀请参加 Deriv 合作伙伴研讨会" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("尊敬的 [合作伙伴姓名], ", $false, $false, $false, $false, $false, $true, 1, $false, "尊敬的 [合作夥伴姓名]， ", 1)
if (-not $ok) { Write-Output "MISS 13: 尊敬的 [合作伙伴姓名], " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("很高兴地通知您，Deriv 联盟团队将于 [MONTH] 在 [CITY] 与您，我们尊贵的合作伙伴会面！", $false, $false, $false, $false, $false, $true, 1, $false, "很高興地通知您，Deriv 聯盟團隊將於 [MONTH] 在 [CITY] 與您，我們尊貴的合作夥伴會面！", 1)
if (-not $ok) { Write-Output "MISS 14: 很高兴地通知您，Deriv 联盟团队将于 [MONTH] 在 [CITY] 与您，我们尊贵的合作伙伴会面！" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("区域经理将于 [DATE] 前通知您确切的地点", $false, $false, $false, $false, $false, $true, 1, $false, "您的區域經理將於 [DATE] 前通知您確切的地點", 1)
if (-not $ok) { Write-Output "MISS 15: 区域经理将于 [DATE] 前通知您确切的地点" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("我们将在为期一天的研讨会提供技术和营销支持，让您与其他合作伙伴共进美味午餐，并听取您对合作伙伴计划的反馈意见。 这是您发表意见的机会，以帮助我们规划未来的工作，为您提供更好的支持。 ", $false, $false, $false, $false, $false, $true, 1, $false, "我們將在為期一天的研討會提供技術和行銷支援，讓您與其他合作夥伴一起享用美味午餐，並聽取您對合作夥伴計劃的反饋意見。 這是您發表意見的機會，以幫助我們計劃未來的工作，為您提供更好的支援。 ", 1)
if (-not $ok) { Write-Output "MISS 16: 我们将在为期一天的研讨会提供技术和营销支持，让您与其他合作伙伴共进美味午餐，并听取您对合作伙伴计划的反馈意见。 这是您发表意见的机会，以帮助我们规划未来的工作，为您提供更好的支持。 " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("作为回复，请在 ", $false, $false, $false, $false, $false, $true, 1, $false, "請在 ", 1)
if (-not $ok) { Write-Output "MISS 17: 作为回复，请在 " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute(" 前提交注册表格。 请注意，参加者以先到先得的方式确认。 期待在那里见到您！", $false, $false, $false, $false, $false, $true, 1, $false, " 前提交註冊表格，確認參加。 請注意，參加者以先到先得的方式確認。 期待在那裡見到您！", 1)
if (-not $ok) { Write-Output "MISS 18:  前提交注册表格。 请注意，参加者以先到先得的方式确认。 期待在那里见到您！" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("发送我的详细信息", $false, $false, $false, $false, $false, $true, 1, $false, "發送我的詳細資訊", 1)
if (-not $ok) { Write-Output "MISS 19: 发送我的详细信息" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("如果有任何疑问，请联系我们：", $false, $false, $false, $false, $false, $true, 1, $false, "如有任何疑問，請與我們聯繫：", 1)
if (-not $ok) { Write-Output "MISS 20: 如果有任何疑问，请联系我们：" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("实时聊天", $false, $false, $false, $false, $false, $true, 1, $false, "即時聊天", 1)
if (-not $ok) { Write-Output "MISS 21: 实时聊天" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("。 / 如有任何疑问，请通过 [电子邮件地址] 或 [WHATSAPP 号码] (WhatsApp) 联系您的区域经理 [NAME] 。 ", $false, $false, $false, $false, $false, $true, 1, $false, "。 / 如有任何疑問，請通過 [電子郵件地址] 或 [WHATSAPP 號碼] (WhatsApp) 聯繫您的區域經理 [NAME] 。 ", 1)
if (-not $ok) { Write-Output "MISS 22: 。 / 如有任何疑问，请通过 [电子邮件地址] 或 [WHATSAPP 号码] (WhatsApp) 联系您的区域经理 [NAME] 。 " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("葡萄牙语", $false, $false, $false, $false, $false, $true, 1, $false, "葡萄牙文", 1)
if (-not $ok) { Write-Output "MISS 23: 葡萄牙语" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("英语", $false, $false, $false, $false, $false, $true, 1, $false, "英文", 1)
if (-not $ok) { Write-Output "MISS 24: 英语" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("主题行", $false, $false, $false, $false, $false, $true, 1, $false, "主題行", 1)
if (-not $ok) { Write-Output "MISS 25: 主题行" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute(": 在 [CITY] 了解我们的团队 | [DATE]", $false, $false, $false, $false, $false, $true, 1, $false, ": 認識我們的團隊在 [CITY] | [DATE]", 1)
if (-not $ok) { Write-Output "MISS 26: : 在 [CITY] 了解我们的团队 | [DATE]" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("您被邀请参加我们的 ", $false, $false, $false, $false, $false, $true, 1, $false, "我們誠摯邀請您參加 ", 1)
if (-not $ok) { Write-Output "MISS 27: 您被邀请参加我们的 " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("Deriv 合作伙伴研讨会", $false, $false, $false, $false, $false, $true, 1, $false, "Deriv 夥伴研討會", 1)
if (-not $ok) { Write-Output "MISS 28: Deriv 合作伙伴研讨会" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("尊敬的 [合作伙伴姓名], ", $false, $false, $false, $false, $false, $true, 1, $false, "親愛的 [合作夥伴姓名]， ", 1)
if (-not $ok) { Write-Output "MISS 29: 尊敬的 [合作伙伴姓名], " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("我们非常高兴地通知您，Deriv 的合作伙伴团队将在 [CITY] 于 [MONTH] 见到您，我们尊贵的合作伙伴！", $false, $false, $false, $false, $false, $true, 1, $false, "我們很高興告訴您，Deriv 的合作夥伴團隊將於 [MONTH] 在 [CITY] 與您會面，我們親愛的夥伴！", 1)
if (-not $ok) { Write-Output "MISS 30: 我们非常高兴地通知您，Deriv 的合作伙伴团队将在 [CITY] 于 [MONTH] 见到您，我们尊贵的合作伙伴！" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("您的国家经理将在 [DATE] 前告知您确切的地点", $false, $false, $false, $false, $false, $true, 1, $false, "您的區域經理將於 [DATE] 前通知您確切的地點", 1)
if (-not $ok) { Write-Output "MISS 31: 您的国家经理将在 [DATE] 前告知您确切的地点" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("在这为期一天的研讨会上，我们将提供技术和市场支持，为您提供与其他合作伙伴社交的机会，并在美味的午餐期间倾听您对我们合作伙伴计划的反馈。 这是您发声的绝佳机会，这将帮助我们规划未来的工作，以便更好地支持您。 ", $false, $false, $false, $false, $false, $true, 1, $false, "在這次為期一天的研討會中，我們將提供技術和市場支持，並提供與其他夥伴互動的機會，享用美味的午餐，同時聆聽您對我們的夥伴計劃的反饋。 這是您讓自己的聲音被聽到的好機會，使我們能夠規劃未來的行動，以便為您提供更好的支持。 ", 1)
if (-not $ok) { Write-Output "MISS 32: 在这为期一天的研讨会上，我们将提供技术和市场支持，为您提供与其他合作伙伴社交的机会，并在美味的午餐期间倾听您对我们合作伙伴计划的反馈。 这是您发声的绝佳机会，这将帮助我们规划未来的工作，以便更好地支持您。 " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("请通过在 ", $false, $false, $false, $false, $false, $true, 1, $false, "請在 ", 1)
if (-not $ok) { Write-Output "MISS 33: 请通过在 " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute(" 前提交注册表格来确认您的出席。 请注意，参加者以先到先得的方式确认。 期待您的到来！ ", $false, $false, $false, $false, $false, $true, 1, $false, " 前提交註冊表格以確認您的參加。 請注意，參加者將以先到先得的方式確認。 期待您的參加！ ", 1)
if (-not $ok) { Write-Output "MISS 34:  前提交注册表格来确认您的出席。 请注意，参加者以先到先得的方式确认。 期待您的到来！ " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("发送我的详细信息", $false, $false, $false, $false, $false, $true, 1, $false, "發送我的詳細資訊", 1)
if (-not $ok) { Write-Output "MISS 35: 发送我的详细信息" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("如有任何疑问，请通过 ", $false, $false, $false, $false, $false, $true, 1, $false, "如有任何疑問，請透過我們的 ", 1)
if (-not $ok) { Write-Output "MISS 36: 如有任何疑问，请通过 " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("实时聊天", $false, $false, $false, $false, $false, $true, 1, $false, "即時聊天", 1)
if (-not $ok) { Write-Output "MISS 37: 实时聊天" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute(" 或者 ", $false, $false, $false, $false, $false, $true, 1, $false, " 或 ", 1)
if (-not $ok) { Write-Output "MISS 38:  或者 " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute(" 联系我们。 / 如有任何疑问，请联系您所在国家的经理 [NAME]，电子邮件 [EMAIL ADDRESS] 或 [WHATSAPP NO] (WhatsApp)。", $false, $false, $false, $false, $false, $true, 1, $false, " 與我們聯繫。 / 如有任何疑問，請與您的區域經理 [NAME] 聯繫，郵件請發送至 [EMAIL ADDRESS] 或 [WHATSAPP NO] (WhatsApp)。", 1)
if (-not $ok) { Write-Output "MISS 39:  联系我们。 / 如有任何疑问，请联系您所在国家的经理 [NAME]，电子邮件 [EMAIL ADDRESS] 或 [WHATSAPP NO] (WhatsApp)。" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("Deriv 合作伙伴团队", $false, $false, $false, $false, $false, $true, 1, $false, "Deriv 合作夥伴團隊", 1)
if (-not $ok) { Write-Output "MISS 40: Deriv 合作伙伴团队" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("法语", $false, $false, $false, $false, $false, $true, 1, $false, "法文", 1)
if (-not $ok) { Write-Output "MISS 41: 法语" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("英语", $false, $false, $false, $false, $false, $true, 1, $false, "英文", 1)
if (-not $ok) { Write-Output "MISS 42: 英语" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("主题行", $false, $false, $false, $false, $false, $true, 1, $false, "主題行", 1)
if (-not $ok) { Write-Output "MISS 43: 主题行" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute(": 在 [CITY] 了解我们的团队 | [DATE]", $false, $false, $false, $false, $false, $true, 1, $false, ": 認識我們的團隊位於 [CITY] | [DATE]", 1)
if (-not $ok) { Write-Output "MISS 44: : 在 [CITY] 了解我们的团队 | [DATE]" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("您被邀请参加我们的 Deriv 研讨会", $false, $false, $false, $false, $false, $true, 1, $false, "您被邀請參加我們的 Deriv 研討會", 1)
if (-not $ok) { Write-Output "MISS 45: 您被邀请参加我们的 Deriv 研讨会" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("亲爱的 [合作伙伴姓名], ", $false, $false, $false, $false, $false, $true, 1, $false, "親愛的 [合作夥伴姓名]， ", 1)
if (-not $ok) { Write-Output "MISS 46: 亲爱的 [合作伙伴姓名], " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("我们很高兴地通知您，Deriv 的合作伙伴团队将在 [CITY] 于 [MONTH] 见到您，我们尊贵的合作伙伴！", $false, $false, $false, $false, $false, $true, 1, $false, "我們非常高興地宣布，Deriv 合作夥伴團隊將於 [MONTH] 在 [CITY] 與您會面，我們珍貴的夥伴！", 1)
if (-not $ok) { Write-Output "MISS 47: 我们很高兴地通知您，Deriv 的合作伙伴团队将在 [CITY] 于 [MONTH] 见到您，我们尊贵的合作伙伴！" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("您的国家经理将在 [DATE] 前告知您确切的地点", $false, $false, $false, $false, $false, $true, 1, $false, "您的國家經理將於 [DATE] 前通知您確切的地點", 1)
if (-not $ok) { Write-Output "MISS 48: 您的国家经理将在 [DATE] 前告知您确切的地点" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("[DATE] | [TIME]", $false, $false, $false, $false, $false, $true, 1, $false, "[DATE] | [HEURE]", 1)
if (-not $ok) { Write-Output "MISS 49: [DATE] | [TIME]" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("在这为期一天的研讨会上，我们将提供技术和市场支持，为您提供与其他合作伙伴社交的机会，并在美味的午餐期间倾听您对我们合作伙伴计划的反馈。 这是您发声的绝佳机会，这将帮助我们规划未来的工作，以便更好地支持您。", $false, $false, $false, $false, $false, $true, 1, $false, "在這次為期一天的研討會中，我們將提供技術和市場支持，並提供與其他夥伴建立聯繫的機會，在豐盛的午餐期間進行交流，並聆聽您的意見，針對我們的夥伴計劃提供反饋。 這是您表達意見的機會，這將幫助我們規劃未來的努力，為您提供更好的支持。", 1)
if (-not $ok) { Write-Output "MISS 50: 在这为期一天的研讨会上，我们将提供技术和市场支持，为您提供与其他合作伙伴社交的机会，并在美味的午餐期间倾听您对我们合作伙伴计划的反馈。 这是您发声的绝佳机会，这将帮助我们规划未来的工作，以便更好地支持您。" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("请在 [DATE] 前提交注册表格确认出席。 请注意，参加者以先到先得的方式确认。 我们期待您的到来！", $false, $false, $false, $false, $false, $true, 1, $false, "請在 [DATE] 前通過提交註冊表格來確認您的出席。 請注意，參加者是以先到先得的方式確認。 我們期待在那裡與您見面！", 1)
if (-not $ok) { Write-Output "MISS 51: 请在 [DATE] 前提交注册表格确认出席。 请注意，参加者以先到先得的方式确认。 我们期待您的到来！" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("发送我的详细信息", $false, $false, $false, $false, $false, $true, 1, $false, "發送我的詳細資訊", 1)
if (-not $ok) { Write-Output "MISS 52: 发送我的详细信息" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("如有任何疑问，请通过实时聊天或 WhatsApp 联系我们。 / 如果您有任何问题，请联系您的国家经理 [NAME]，如 [EMAIL ADDRESS] 或 [WHATSAPP NO] (WhatsApp)。", $false, $false, $false, $false, $false, $true, 1, $false, "Si vous avez des questions, veuillez nous contacter via le chat en direct ou WhatsApp sur notre site Web. / Si vous avez des questions, veuillez contacter votre gestionnaire de compte [NOM] sur [ADRESSE EMAIL] ou [WHATSAPP NO] (WhatsApp).(", 1)
if (-not $ok) { Write-Output "MISS 53: 如有任何疑问，请通过实时聊天或 WhatsApp 联系我们。 / 如果您有任何问题，请联系您的国家经理 [NAME]，如 [EMAIL ADDRESS] 或 [WHATSAPP NO] (WhatsApp)。" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("Deriv 合作伙伴团队", $false, $false, $false, $false, $false, $true, 1, $false, "Deriv 合作夥伴團隊", 1)
if (-not $ok) { Write-Output "MISS 54: Deriv 合作伙伴团队" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("泰语", $false, $false, $false, $false, $false, $true, 1, $false, "泰文", 1)
if (-not $ok) { Write-Output "MISS 55: 泰语" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("英语", $false, $false, $false, $false, $false, $true, 1, $false, "英文", 1)
if (-not $ok) { Write-Output "MISS 56: 英语" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("主题行", $false, $false, $false, $false, $false, $true, 1, $false, "主題行", 1)
if (-not $ok) { Write-Output "MISS 57: 主题行" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute(": 在 [CITY] 了解我们的团队 | [DATE]", $false, $false, $false, $false, $false, $true, 1, $false, ": 在 [CITY] | [DATE] 與我們的團隊會面", 1)
if (-not $ok) { Write-Output "MISS 58: : 在 [CITY] 了解我们的团队 | [DATE]" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("您被邀请参加我们的 ", $false, $false, $false, $false, $false, $true, 1, $false, "您被邀請參加我們的研討會 ", 1)
if (-not $ok) { Write-Output "MISS 59: 您被邀请参加我们的 " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("Deriv 合作伙伴研讨会", $false, $false, $false, $false, $false, $true, 1, $false, "我們的 Deriv 夥伴研討會", 1)
if (-not $ok) { Write-Output "MISS 60: Deriv 合作伙伴研讨会" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("亲爱的 [合作伙伴姓名] ", $false, $false, $false, $false, $false, $true, 1, $false, "親愛的 [合作夥伴姓名] ", 1)
if (-not $ok) { Write-Output "MISS 61: 亲爱的 [合作伙伴姓名] " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("我们非常高兴地通知您，Deriv 的合作伙伴团队将在 [CITY] 于 [MONTH] 见到您，我们尊贵的合作伙伴！", $false, $false, $false, $false, $false, $true, 1, $false, "我們很高興通知您，Deriv 的合作夥伴團隊將於 [MONTH] 在 [CITY] 與您會面，您是我們珍貴的夥伴！", 1)
if (-not $ok) { Write-Output "MISS 62: 我们非常高兴地通知您，Deriv 的合作伙伴团队将在 [CITY] 于 [MONTH] 见到您，我们尊贵的合作伙伴！" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("您的国家经理将在 [DATE] 前告知您确切的地点", $false, $false, $false, $false, $false, $true, 1, $false, "您的國家經理將在 [DATE] 前告訴您精確的位置", 1)
if (-not $ok) { Write-Output "MISS 63: 您的国家经理将在 [DATE] 前告知您确切的地点" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("在这一天的研讨会上，我们将提供技术和市场支持，为您提供与其他合作伙伴交流的机会，并倾听您关于我们合作伙伴计划的反馈。 ", $false, $false, $false, $false, $false, $true, 1, $false, "在這一天的研討會中，我們將提供技術和市場支持，提供與其他夥伴建立聯繫的機會，在美味的午餐上交流，並聽取您對我們夥伴計劃的反饋。", 1)
if (-not $ok) { Write-Output "MISS 64: 在这一天的研讨会上，我们将提供技术和市场支持，为您提供与其他合作伙伴交流的机会，并倾听您关于我们合作伙伴计划的反馈。 " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("โปรดตอบรับคำเชิญโดยการส่งแบบฟอร์มลงทะเบียนภายในวันที่ ", $false, $false, $false, $false, $false, $true, 1, $false, "請在 ", 1)
if (-not $ok) { Write-Output "MISS 65: โปรดตอบรับคำเชิญโดยการส่งแบบฟอร์มลงทะเบียนภายในวันที่ " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute(" โปรดทราบว่า การเข้าร่วมกิจกรรมจะได้รับการยืนยันตาม", $false, $false, $false, $false, $false, $true, 1, $false, " 前提交註冊表，以確認您的參加。請注意，參加者的確認將依照", 1)
if (-not $ok) { Write-Output "MISS 66:  โปรดทราบว่า การเข้าร่วมกิจกรรมจะได้รับการยืนยันตาม" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("หลัก มาก่อนได้ก่อน หากมีการจองซ้ำกัน", $false, $false, $false, $false, $false, $true, 1, $false, "先到先得", 1)
if (-not $ok) { Write-Output "MISS 67: หลัก มาก่อนได้ก่อน หากมีการจองซ้ำกัน" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute(" และเราหวังว่าจะได้พบคุณที่นั่น!", $false, $false, $false, $false, $false, $true, 1, $false, "原則。期待在研討會上見到您！", 1)
if (-not $ok) { Write-Output "MISS 68:  และเราหวังว่าจะได้พบคุณที่นั่น!" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("发送我的详细信息", $false, $false, $false, $false, $false, $true, 1, $false, "發送我的資料", 1)
if (-not $ok) { Write-Output "MISS 69: 发送我的详细信息" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("如有疑问，请通过 ", $false, $false, $false, $false, $false, $true, 1, $false, "如您有任何疑問，請通過我們的 ", 1)
if (-not $ok) { Write-Output "MISS 70: 如有疑问，请通过 " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("实时聊天", $false, $false, $false, $false, $false, $true, 1, $false, "即時聊天", 1)
if (-not $ok) { Write-Output "MISS 71: 实时聊天" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute(" 联系我们。如有任何问题，请联系您国家的经理 [NAME]，电子邮件 [EMAIL ADDRESS] 或 [WHATSAPP NO] (WhatsApp)。", $false, $false, $false, $false, $false, $true, 1, $false, " 聯繫我們。如果您有任何問題，請聯繫您所在國的經理 [NAME]，或通過 [EMAIL ADDRESS] 或 [WHATSAPP NO] (WhatsApp)。", 1)
if (-not $ok) { Write-Output "MISS 72:  联系我们。如有任何问题，请联系您国家的经理 [NAME]，电子邮件 [EMAIL ADDRESS] 或 [WHATSAPP NO] (WhatsApp)。" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("来自 Deriv 合作伙伴团队", $false, $false, $false, $false, $false, $true, 1, $false, "來自 Deriv 合作夥伴團隊", 1)
if (-not $ok) { Write-Output "MISS 73: 来自 Deriv 合作伙伴团队" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("越南语", $false, $false, $false, $false, $false, $true, 1, $false, "越南文", 1)
if (-not $ok) { Write-Output "MISS 74: 越南语" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("英语", $false, $false, $false, $false, $false, $true, 1, $false, "英文", 1)
if (-not $ok) { Write-Output "MISS 75: 英语" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("主题行", $false, $false, $false, $false, $false, $true, 1, $false, "主題行", 1)
if (-not $ok) { Write-Output "MISS 76: 主题行" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute(": 期待在 [CITY] 与您相见 | [DATE]", $false, $false, $false, $false, $false, $true, 1, $false, ": Deriv 期待在 [CITY] | [DATE] 與您見面", 1)
if (-not $ok) { Write-Output "MISS 77: : 期待在 [CITY] 与您相见 | [DATE]" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("诚邀您参加 Deriv 合作伙伴研讨会", $false, $false, $false, $false, $false, $true, 1, $false, "我們誠摯邀請您參加 Deriv 夥伴研討會", 1)
if (-not $ok) { Write-Output "MISS 78: 诚邀您参加 Deriv 合作伙伴研讨会" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("您好 [合作伙伴姓名], ", $false, $false, $false, $false, $false, $true, 1, $false, "您好 [合作夥伴姓名]， ", 1)
if (-not $ok) { Write-Output "MISS 79: 您好 [合作伙伴姓名], " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("我们很高兴地通知您，Deriv 的合作伙伴团队将在 [CITY] 的 [MONTH] 见面，与您这位重要的合作伙伴见面！", $false, $false, $false, $false, $false, $true, 1, $false, "我們很高興通知您，Deriv 的合作夥伴團隊將於 [MONTH] 在 [CITY] 與您會面，您是我們重要的夥伴！", 1)
if (-not $ok) { Write-Output "MISS 80: 我们很高兴地通知您，Deriv 的合作伙伴团队将在 [CITY] 的 [MONTH] 见面，与您这位重要的合作伙伴见面！" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("您的国家经理将在 [DATE] 前告知您确切的地点", $false, $false, $false, $false, $false, $true, 1, $false, "您的國家經理將在 [DATE] 前通知您確切的地點", 1)
if (-not $ok) { Write-Output "MISS 81: 您的国家经理将在 [DATE] 前告知您确切的地点" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("在为期一天的研讨会上，我们将提供技术和市场支持，为您提供与其他合作伙伴的交流机会，并倾听您对我们合作伙伴计划的反馈。 这是您分享意见的机会，这将帮助我们计划未来的工作，以便更好地支持您。 ", $false, $false, $false, $false, $false, $true, 1, $false, "在為期一天的研討會中，我們將提供技術和市場的支持，讓您有機會在午餐中與其他夥伴建立聯繫，並聆聽您對我們夥伴計劃的反饋。 這是讓您分享意見的機會，這將幫助我們計劃更好的未來支持。 ", 1)
if (-not $ok) { Write-Output "MISS 82: 在为期一天的研讨会上，我们将提供技术和市场支持，为您提供与其他合作伙伴的交流机会，并倾听您对我们合作伙伴计划的反馈。 这是您分享意见的机会，这将帮助我们计划未来的工作，以便更好地支持您。 " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("请在 ", $false, $false, $false, $false, $false, $true, 1, $false, "請通過在 ", 1)
if (-not $ok) { Write-Output "MISS 83: 请在 " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute(" 前回复我们，提交注册表单。 请注意，先到先得，先确认参与。 期待在研讨会上见到您！", $false, $false, $false, $false, $false, $true, 1, $false, " 前提交申請表來確認您的出席。 請注意，我們將優先確認較早提交的申請。 期待在研討會上見到您！", 1)
if (-not $ok) { Write-Output "MISS 84:  前回复我们，提交注册表单。 请注意，先到先得，先确认参与。 期待在研讨会上见到您！" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("发送我的信息", $false, $false, $false, $false, $false, $true, 1, $false, "發送我的詳細資料", 1)
if (-not $ok) { Write-Output "MISS 85: 发送我的信息" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("如果您需要支持，请通过 ", $false, $false, $false, $false, $false, $true, 1, $false, "如果您需要支持，請通過網站上的 ", 1)
if (-not $ok) { Write-Output "MISS 86: 如果您需要支持，请通过 " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("实时聊天", $false, $false, $false, $false, $false, $true, 1, $false, "即時聊天", 1)
if (-not $ok) { Write-Output "MISS 87: 实时聊天" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute(" 联系我们。 如果您有任何问题，请联系您的国家经理 [NAME]，通过 [EMAIL ADDRESS] 或 [WHATSAPP NO] (WhatsApp)。 ", $false, $false, $false, $false, $false, $true, 1, $false, " 與我們聯繫。 如果您有任何疑問，請聯繫您的國家經理 [NAME]，通過 [EMAIL ADDRESS] 或 [WHATSAPP NO] (WhatsApp)。 ", 1)
if (-not $ok) { Write-Output "MISS 88:  联系我们。 如果您有任何问题，请联系您的国家经理 [NAME]，通过 [EMAIL ADDRESS] 或 [WHATSAPP NO] (WhatsApp)。 " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("Deriv 合作伙伴团队", $false, $false, $false, $false, $false, $true, 1, $false, "Deriv 合作夥伴團隊", 1)
if (-not $ok) { Write-Output "MISS 89: Deriv 合作伙伴团队" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("西班牙语", $false, $false, $false, $false, $false, $true, 1, $false, "西班牙文", 1)
if (-not $ok) { Write-Output "MISS 90: 西班牙语" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("英语", $false, $false, $false, $false, $false, $true, 1, $false, "英文", 1)
if (-not $ok) { Write-Output "MISS 91: 英语" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("主题行", $false, $false, $false, $false, $false, $true, 1, $false, "主題行", 1)
if (-not $ok) { Write-Output "MISS 92: 主题行" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute(": 在 [CITY] 了解我们的团队 | [DATE]", $false, $false, $false, $false, $false, $true, 1, $false, ": 認識我們的團隊在 [CITY] | [DATE]", 1)
if (-not $ok) { Write-Output "MISS 93: : 在 [CITY] 了解我们的团队 | [DATE]" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("我们诚邀您参加 Deriv 为合作伙伴举办的研讨会", $false, $false, $false, $false, $false, $true, 1, $false, "我們誠摯邀請您參加 Deriv 夥伴研討會", 1)
if (-not $ok) { Write-Output "MISS 94: 我们诚邀您参加 Deriv 为合作伙伴举办的研讨会" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("亲爱的 [合作伙伴姓名], ", $false, $false, $false, $false, $false, $true, 1, $false, "尊敬的 [合作夥伴姓名]， ", 1)
if (-not $ok) { Write-Output "MISS 95: 亲爱的 [合作伙伴姓名], " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("我们很高兴地通知您，Deriv 的合作伙伴团队将在 [CITY] 于 [MONTH] 见到您，我们尊贵的合作伙伴！", $false, $false, $false, $false, $false, $true, 1, $false, "我們很高興通知您，Deriv 的合作夥伴團隊將於 [MONTH] 在 [CITY] 與您會面，您是我們珍貴的夥伴！", 1)
if (-not $ok) { Write-Output "MISS 96: 我们很高兴地通知您，Deriv 的合作伙伴团队将在 [CITY] 于 [MONTH] 见到您，我们尊贵的合作伙伴！" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("您的经理将于 [DATE] 前告知您确切位置", $false, $false, $false, $false, $false, $true, 1, $false, "您的經理將於 [DATE] 通知您確切的位置", 1)
if (-not $ok) { Write-Output "MISS 97: 您的经理将于 [DATE] 前告知您确切位置" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("在这为期一天的研讨会上，我们将提供技术和市场支持，为您提供与其他合作伙伴互动的机会，并倾听您对我们合作伙伴计划的反馈。 这是您发声的机会，这将帮助我们计划未来的努力，能够更好地支持您。 ", $false, $false, $false, $false, $false, $true, 1, $false, "在這一天的研討會中，我們將提供技術和市場支持，並提供與其他夥伴建立聯繫的機會，享用美味的午餐，同時聆聽您對我們夥伴計劃的反饋。 這是讓您發出聲音的機會，這將有助於我們規劃未來的努力，以更好地支持您。 ", 1)
if (-not $ok) { Write-Output "MISS 98: 在这为期一天的研讨会上，我们将提供技术和市场支持，为您提供与其他合作伙伴互动的机会，并倾听您对我们合作伙伴计划的反馈。 这是您发声的机会，这将帮助我们计划未来的努力，能够更好地支持您。 " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("请在 ", $false, $false, $false, $false, $false, $true, 1, $false, "請在 ", 1)
if (-not $ok) { Write-Output "MISS 99: 请在 " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute(" 前提交注册表确认您的出席。 请注意，出席将根据先到先得的原则确认。 期待您的光临。", $false, $false, $false, $false, $false, $true, 1, $false, " 前提交註冊表以確認您的參加。 請注意，參加者的確認將依照先到先得的原則。 期待與您會面。", 1)
if (-not $ok) { Write-Output "MISS 100:  前提交注册表确认您的出席。 请注意，出席将根据先到先得的原则确认。 期待您的光临。" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("发送我的详细信息", $false, $false, $false, $false, $false, $true, 1, $false, "發送我的詳細資料", 1)
if (-not $ok) { Write-Output "MISS 101: 发送我的详细信息" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("如有任何问题，请通过 ", $false, $false, $false, $false, $false, $true, 1, $false, "如果您有任何問題，請通過網站上的 ", 1)
if (-not $ok) { Write-Output "MISS 102: 如有任何问题，请通过 " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("实时聊天", $false, $false, $false, $false, $false, $true, 1, $false, "即時聊天", 1)
if (-not $ok) { Write-Output "MISS 103: 实时聊天" } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute(" 联系我们。 / 如果您有任何问题，请联系您国家的经理 [NAME]，在 [EMAIL ADDRESS] 或 [WHATSAPP NO] (WhatsApp) 联系。 ", $false, $false, $false, $false, $false, $true, 1, $false, " 聯繫我們。 / 如果您有任何問題，請聯繫您所在國的經理 [NAME]，郵件請發送至 [EMAIL ADDRESS] 或 [WHATSAPP NO] (WhatsApp)。 ", 1)
if (-not $ok) { Write-Output "MISS 104:  联系我们。 / 如果您有任何问题，请联系您国家的经理 [NAME]，在 [EMAIL ADDRESS] 或 [WHATSAPP NO] (WhatsApp) 联系。 " } else { $count = $count + 1 }

$rng = $d.Content
$ok = $rng.Find.Execute("Deriv 合作伙伴团队", $false, $false, $false, $false, $false, $true, 1, $false, "Deriv 合作夥伴團隊", 1)
if (-not $ok) { Write-Output "MISS 105: Deriv 合作伙伴团队" } else { $count = $count + 1 }

Write-Output "Replaced: $count of 106"

# Special case: French support paragraph gains a new trailing ")" run
# after the WhatsApp hyperlink (inside the same paragraph, before the paragraph mark).
$rng3 = $d.Content
$foundMarker = $rng3.Find.Execute("Si vous avez des questions, veuillez nous contacter via le chat en direct", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundMarker) {
  $para = $rng3.Paragraphs(1)
  $prng = $para.Range
  $endRng = $prng.Duplicate
  [void]$endRng.MoveEnd(1, -1)
  $endRng.Collapse(0)
  $endRng.InsertAfter(")")
  Write-Output "Inserted closing paren"
} else {
  Write-Output "MISS marker for closing paren insertion"
}
